# Saldo.xlsx update — refresh the "Export" sheet's top balance rows with a
# newer export: two of the former rows (PHYLIA, MARCELO) are replaced by two
# rows each (net +1 row overall), RODRIGO's balance is updated, and CESAR's
# row is dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: row 5 (currently "004690692 PHYLIA 15273.35") shifts down,
# along with everything below it, opening a fresh blank row at 5. That nets
# the table one extra row, which is exactly what's needed once rows 5-9 are
# rewritten as the six rows below (5-10).
$ws.Rows.Item(5).Insert()

# Column A holds account numbers with significant leading zeros, so force
# text formatting before writing them (otherwise Excel's General format
# auto-converts the numeric-looking string and drops the leading zeros).
$ws.Range("A5:A10").NumberFormat = "@"

$ws.Cells.Item(5,1).Value = "004267119"
$ws.Cells.Item(5,2).Value = "ANA"
$ws.Cells.Item(5,3).Value = 14593.13

$ws.Cells.Item(6,1).Value = "005064906"
$ws.Cells.Item(6,2).Value = "BERTILLA"
$ws.Cells.Item(6,3).Value = 9984.8

$ws.Cells.Item(7,1).Value = "004479734"
$ws.Cells.Item(7,2).Value = "RODRIGO"
$ws.Cells.Item(7,3).Value = 7165.3

$ws.Cells.Item(8,1).Value = "005624730"
$ws.Cells.Item(8,2).Value = "ISABEL"
$ws.Cells.Item(8,3).Value = 5000

$ws.Cells.Item(9,1).Value = "004260002"
$ws.Cells.Item(9,2).Value = "ERICA"
$ws.Cells.Item(9,3).Value = 4279.88

$ws.Cells.Item(10,1).Value = "004267044"
$ws.Cells.Item(10,2).Value = "PATRICIA"
$ws.Cells.Item(10,3).Value = 2490.72
